$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously held, in column A:
#   A7 = "data"          (style 0)
#   A8 = 1                (style 1)
#   A9 = "xsd:decimal"    (style 1)
# It now needs to hold:
#   A7 = 1                (style 0, unchanged from before)
#   A8 = "xsd:decimal"    (style 1, unchanged from before)
#   A9 = "data"           (style 0 -- picks up the look A7 used to have)

# Give A9 the same formatting A7 currently has (both are visually identical
# "General" Arial styles, but this keeps the underlying style index in sync
# with the target file) before the content shuffles around.
$ws.Range("A7").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats

# Now move the values into their new homes.
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = "xsd:decimal"
$ws.Range("A9").Value = "data"

# Update the active selection to match: single cell A8 selected.
$ws.Range("A8").Select()
